$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Task Schedule")
$ws.Activate()

# --- Delete row 44 ("0. Technology Test: Vuforia") ---------------------
# This shifts rows 45..986 up by one (new last row = 985) and also drops
# the now-unused "Vuforia" shared string automatically on save.
$ws.Rows(44).Delete()

# --- Fill in the previously-empty State (J) cells for rows 20-23 -------
$ws.Range("J20").Value = "Done"
$ws.Range("J21").Value = "Done"
$ws.Range("J22").Value = "Ongoing"
$ws.Range("J23").Value = "Ongoing"

# --- Clear the State (J) cells for rows 26-29 (was "Delayed") ----------
$ws.Range("J26").ClearContents()
$ws.Range("J27").ClearContents()
$ws.Range("J28").ClearContents()
$ws.Range("J29").ClearContents()

# --- Hide columns D and G ------------------------------------------------
$ws.Columns("D").Hidden = $true
$ws.Columns("G").Hidden = $true

# --- Update the data-validation list source (now J39:J45, was J39:J46) -
# Drop it from the now-gone single cell J44 leftover reference and from
# J2:J38, then re-apply cleanly on J2:J38 with the corrected source range.
$ws.Range("J2:J38").Validation.Delete()
$ws.Range("J2:J38").Validation.Add(3, 1, 1, "=`$J`$39:`$J`$45")
$ws.Range("J2:J38").Validation.IgnoreBlank = $true
$ws.Range("J2:J38").Validation.InCellDropdown = $true
$ws.Range("J2:J38").Validation.ShowInput = $true
$ws.Range("J2:J38").Validation.ShowError = $true

# --- Re-anchor the conditional format rule to J2:J38 only (drop J44) ---
$cf = $ws.Range("J2:J38").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("J2:J38"))

# --- Restore the view: scrolled to row 13, active cell L24 -------------
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L24").Select()
